$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 3-17 entirely (data trimmed down to a single data row)
$ws.Range("A3:B17").EntireRow.Delete()

# Update the remaining data row (row 2) with the new date/value pair
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 5.361718827437545
